# Search Remit Test Cases
# Adds new SQL test-case rows (44-49) to the "SQL" sheet and fixes the
# hard-coded TIN in the existing "Sno 43" query to use the {$tin} token.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SQL")

# ---------------------------------------------------------------------
# 1) Fix the existing row 44 query (Sno 43 - subscriber id lookup):
#    hard-coded TIN '840611484' -> templated '{$tin}' (also drops the
#    stray space before the '=').
# ---------------------------------------------------------------------
$ws.Range("B44").Value = @'
Select sr.SBSCR_ID, cp.SETL_DT
from PP001.CONSOLIDATED_PAYMENT cp, PP001.PROVIDER p, PP001.CLAIM c, PP001.UNCONSOLIDATED_PAYMENT ucp, PP001.CLAIM_UNCONSOLIDATED_PAYMENT cup, PP001.SUBSCRIBER sr
where cp.prov_key_id = p.prov_key_id
and cp.CONSL_PAY_NBR = ucp.CONSL_PAY_NBR
and ucp.UCONSL_PAY_KEY_ID = cup.UCONSL_PAY_KEY_ID
and cup.CLM_KEY_ID = c.CLM_KEY_ID
and c.SBSCR_KEY_ID = sr.SBSCR_KEY_ID
and p.PROV_TAX_ID_NBR ='{$tin}'
and cp.SETL_DT <= current date 
order by cp.SETL_DT DESC
fetch first row only
'@

# ---------------------------------------------------------------------
# 2) New row 45 (Sno 44 - provider NPI lookup)
# ---------------------------------------------------------------------
$ws.Range("A44:B44").Copy($ws.Range("A45:B45")) | Out-Null
$ws.Range("A45").Formula = "'44"
$ws.Range("B45").Value = @'
Select p.PROV_NPI_NBR, cp.SETL_DT
from PP001.CONSOLIDATED_PAYMENT cp, PP001.PROVIDER p
where cp.prov_key_id = p.prov_key_id
and p.PROV_TAX_ID_NBR ='{$tin}'
and cp.SETL_DT <= current date 
order by cp.SETL_DT DESC
fetch first row only
'@
$ws.Rows.Item(45).RowHeight = 105

# ---------------------------------------------------------------------
# 3) New row 46 (Sno 45 - claim number lookup)
# ---------------------------------------------------------------------
$ws.Range("A44:B44").Copy($ws.Range("A46:B46")) | Out-Null
$ws.Range("A46").Formula = "'45"
$ws.Range("B46").Value = @'
Select c.CLM_NBR, cp.SETL_DT
from PP001.CONSOLIDATED_PAYMENT cp, PP001.PROVIDER p, PP001.CLAIM c, PP001.UNCONSOLIDATED_PAYMENT ucp, PP001.CLAIM_UNCONSOLIDATED_PAYMENT cup
where cp.prov_key_id = p.prov_key_id
and cp.CONSL_PAY_NBR = ucp.CONSL_PAY_NBR
and ucp.UCONSL_PAY_KEY_ID = cup.UCONSL_PAY_KEY_ID
and cup.CLM_KEY_ID = c.CLM_KEY_ID
and p.PROV_TAX_ID_NBR = '{$tin}'
and cp.SETL_DT <= current date 
order by cp.SETL_DT DESC
fetch first row only
'@
$ws.Rows.Item(46).RowHeight = 165

# ---------------------------------------------------------------------
# 4) New row 47 (Sno 46 - patient first/last name lookup)
# ---------------------------------------------------------------------
$ws.Range("A44:B44").Copy($ws.Range("A47:B47")) | Out-Null
$ws.Range("A47").Formula = "'46"
$ws.Range("B47").Value = @'
Select c.PTNT_FST_NM, c.PTNT_LST_NM, cp.SETL_DT
from PP001.CONSOLIDATED_PAYMENT cp, PP001.PROVIDER p, PP001.CLAIM c, PP001.UNCONSOLIDATED_PAYMENT ucp, PP001.CLAIM_UNCONSOLIDATED_PAYMENT cup
where cp.prov_key_id = p.prov_key_id
and cp.CONSL_PAY_NBR = ucp.CONSL_PAY_NBR
and ucp.UCONSL_PAY_KEY_ID = cup.UCONSL_PAY_KEY_ID
and cup.CLM_KEY_ID = c.CLM_KEY_ID
and p.PROV_TAX_ID_NBR = '{$tin}'
and cp.SETL_DT <= current date 
order by cp.SETL_DT DESC
fetch first row only
'@
$ws.Rows.Item(47).RowHeight = 165

# ---------------------------------------------------------------------
# 5) New row 48 (Sno 47 - zero dollar payment TIN lookup)
# ---------------------------------------------------------------------
$ws.Range("A44:B44").Copy($ws.Range("A48:B48")) | Out-Null
$ws.Range("A48").Formula = "'47"
$ws.Range("B48").Value = @'
Select p.PROV_TAX_ID_NBR
from PP001.PROVIDER p, PP001.UNCONSOLIDATED_PAYMENT ucp, PP001.CLAIM_UNCONSOLIDATED_PAYMENT cup
where ucp.prov_key_id = p.prov_key_id
and ucp.UCONSL_PAY_KEY_ID = cup.UCONSL_PAY_KEY_ID
and cup.CLM_PAY_AMT = '0.00'
and date(ucp.PROC_DTTM) between (current date - 30 days) and (current date - 5 days)
fetch first row only
'@
$ws.Rows.Item(48).RowHeight = 120

# ---------------------------------------------------------------------
# 6) New row 49 (Sno 48 - no query yet, column A only)
# ---------------------------------------------------------------------
$ws.Range("A44").Copy($ws.Range("A49")) | Out-Null
$ws.Range("A49").Formula = "'48"

# ---------------------------------------------------------------------
# 7) Move the selection to where the author left off editing
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("B48").Select() | Out-Null
